$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.Style = "Normal"
}

Set-TextValue "D2" "60.740.03"
Set-TextValue "E2" "  -3.17%  "
Set-TextValue "D3" "2.899.40"
Set-TextValue "E3" "  -4.21%  "
Set-TextValue "E4" "  +0.09%  "
Set-TextValue "D5" "585.23"
Set-TextValue "E5" "  -1.90%  "
Set-TextValue "D6" "145.04"
Set-TextValue "E6" "  -4.22%  "
Set-TextValue "E7" "  +0.00%  "
Set-TextValue "D8" "0.501"
Set-TextValue "E8" "  -3.47%  "
Set-TextValue "D9" "2.897.63"
Set-TextValue "E9" "  -4.16%  "
Set-TextValue "D10" "6.69"
Set-TextValue "E10" "  +4.49%  "
Set-TextValue "E11" "  -5.19%  "
Set-TextValue "D12" "0.445"
Set-TextValue "E12" "  -3.24%  "
Set-TextValue "D13" "0.0000223"
Set-TextValue "E13" "  -4.78%  "
Set-TextValue "D14" "33.36"
Set-TextValue "E14" "  -4.09%  "
Set-TextValue "E15" "  +0.55%  "
Set-TextValue "D16" "3.380.45"
Set-TextValue "E16" "  -4.20%  "
Set-TextValue "D17" "60.660.60"
Set-TextValue "E17" "  -3.20%  "
Set-TextValue "D18" "6.72"
Set-TextValue "E18" "  -4.24%  "
Set-TextValue "D19" "2.901.08"
Set-TextValue "E19" "  -4.21%  "
Set-TextValue "D20" "423.55"
Set-TextValue "E20" "  -5.92%  "
Set-TextValue "D21" "13.51"
Set-TextValue "E21" "  -4.84%  "
Set-TextValue "D22" "0.668"
Set-TextValue "E22" "  -3.41%  "
Set-TextValue "D23" "7.06"
Set-TextValue "E23" "  -5.45%  "
Set-TextValue "D24" "79.83"
Set-TextValue "E24" "  -3.08%  "
Set-TextValue "D25" "11.00"
Set-TextValue "E25" "  +0.87%  "
Set-TextValue "D26" "2.20"
Set-TextValue "E26" "  -2.92%  "
Set-TextValue "D27" "11.80"
Set-TextValue "E27" "  -3.01%  "
Set-TextValue "E28" "  -0.08%  "
Set-TextValue "D29" "1.00"
Set-TextValue "E29" "  +0.17%  "
Set-TextValue "D30" "7.15"
Set-TextValue "E30" "  -3.16%  "
Set-TextValue "D31" "2.61"
Set-TextValue "E31" "  -4.04%  "
Set-TextValue "D32" "2.15"
Set-TextValue "E32" "  -0.13%  "
Set-TextValue "D33" "26.22"
Set-TextValue "E33" "  -4.98%  "
Set-TextValue "E34" "  -4.13%  "
Set-TextValue "D35" "0.0₃0833"
Set-TextValue "E35" "  -3.35%  "
Set-TextValue "D36" "1.00"
Set-TextValue "E36" "  -2.88%  "
Set-TextValue "D37" "5.58"
Set-TextValue "E37" "  -4.96%  "
Set-TextValue "D38" "49.34"
Set-TextValue "E38" "  -1.96%  "
Set-TextValue "D39" "2.95"
Set-TextValue "E39" "  -4.89%  "
Set-TextValue "D40" "0.124"
Set-TextValue "E40" "  +0.01%  "
Set-TextValue "D41" "2.01"
Set-TextValue "E41" "  -4.27%  "
Set-TextValue "D42" "8.64"
Set-TextValue "E42" "  -4.40%  "
Set-TextValue "D43" "0.289"
Set-TextValue "E43" "  -2.58%  "
Set-TextValue "D44" "41.04"
Set-TextValue "E44" "  -0.08%  "
Set-TextValue "D45" "0.0345"
Set-TextValue "E45" "  -3.09%  "
Set-TextValue "D46" "371.43"
Set-TextValue "E46" "  -5.44%  "
Set-TextValue "D47" "2.657.99"
Set-TextValue "E47" "  -2.95%  "
Set-TextValue "D48" "132.96"
Set-TextValue "E48" "  +0.31%  "
Set-TextValue "E49" "  -0.01%  "
Set-TextValue "D50" "25.12"
Set-TextValue "E50" "  +4.76%  "
Set-TextValue "E51" "  -2.07%  "
